$d = $word.ActiveDocument

# The document starts with a paragraph containing "." followed by a long run
# of empty filler paragraphs before the two picture paragraphs. The edit
# removes the "." paragraph and all but one of the filler paragraphs
# (keeping exactly one empty paragraph directly before the first picture).
$start = $d.Paragraphs(1).Range.Start
$end = $d.Paragraphs(61).Range.End
$fillerRange = $d.Range($start, $end)
$fillerRange.Delete()

# Paragraph 2 is now the first picture paragraph. Its run currently begins
# with <w:lastRenderedPageBreak/> before the <w:drawing>; round-tripping the
# range through WordOpenXML/InsertXML drops that stale rendering artifact,
# and then we mark the run NoProof (<w:rPr><w:noProof/></w:rPr>) same as Word
# does when it inserts a picture via "Paste"/"Insert Picture".
$pic1 = $d.Paragraphs(2)
$pic1Range = $pic1.Range
$pic1Xml = $pic1Range.WordOpenXML
$pic1Range.InsertXML($pic1Xml)
$d.Paragraphs(2).Range.NoProofing = 1

# Paragraph 3 is the second picture paragraph; it only needs the NoProof
# run property added.
$pic2 = $d.Paragraphs(3)
$pic2.Range.NoProofing = 1
